$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.016.13"
$ws.Range("E2").Value = "  +0.73%  "

$ws.Range("D3").Value = "2.466.84"
$ws.Range("E3").Value = "  +0.71%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "560.06"
$ws.Range("E5").Value = "  -0.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.16"
$ws.Range("E6").Value = "  -0.43%  "

$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.505"
$ws.Range("E8").Value = "  +0.08%  "

$ws.Range("D9").Value = "2.465.49"
$ws.Range("E9").Value = "  +0.70%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.150"
$ws.Range("E10").Value = "  -0.67%  "

$ws.Range("E11").Value = "  +0.67%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.87"
$ws.Range("E12").Value = "  +1.02%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.331"
$ws.Range("E13").Value = "  -3.32%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "2.905.82"
$ws.Range("E14").Value = "  +0.11%  "

$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "68.890.71"
$ws.Range("E15").Value = "  +0.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000168"
$ws.Range("E16").Value = "  -1.63%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "23.58"
$ws.Range("E17").Value = "  -0.49%  "

$ws.Range("D18").Value = "2.469.93"
$ws.Range("E18").Value = "  +0.52%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.72"
$ws.Range("E19").Value = "  -1.73%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "336.84"
$ws.Range("E20").Value = "  -3.37%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.97"
$ws.Range("E21").Value = "  -3.13%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.79"
$ws.Range("E22").Value = "  -0.69%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.88"
$ws.Range("E23").Value = "  -0.06%  "

$ws.Range("E24").Value = "  -0.08%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "66.97"
$ws.Range("E25").Value = "  -1.90%  "

$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "2.600.84"
$ws.Range("E26").Value = "  +0.69%  "

$ws.Range("B27").Value = "NEARProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.66"
$ws.Range("E27").Value = "  -2.77%  "

$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.998"
$ws.Range("E28").Value = "  -1.68%  "

$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.22"
$ws.Range("E29").Value = "  -0.80%  "

$ws.Range("D30").Value = "0.0₃0817"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.19"
$ws.Range("E31").Value = "  -1.91%  "

$ws.Range("E32").Value = "  -0.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "429.73"
$ws.Range("E33").Value = "  -0.17%  "

$ws.Range("E34").Value = "  -3.30%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.62"
$ws.Range("E35").Value = "  -3.62%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "156.86"
$ws.Range("E36").Value = "  -0.29%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.05"
$ws.Range("E37").Value = "  +0.27%  "

$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("E39").Value = "  -1.10%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.78"
$ws.Range("E40").Value = "  -0.86%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.300"
$ws.Range("E41").Value = "  -2.33%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.41"
$ws.Range("E42").Value = "  -2.33%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.46"
$ws.Range("E43").Value = "  -4.83%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.07"
$ws.Range("E44").Value = "  -2.15%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.06"
$ws.Range("E45").Value = "  -1.56%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "132.30"
$ws.Range("E46").Value = "  -2.25%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.35"
$ws.Range("E47").Value = "  -0.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0715"
$ws.Range("E48").Value = "  -0.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.483"
$ws.Range("E49").Value = "  -1.37%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.560"
$ws.Range("E50").Value = "  -0.64%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0914"
$ws.Range("E51").Value = "  -0.10%  "
